$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Delete empty paragraph before "PUBLICATIONS" (para 73)
# ---------------------------------------------------------------
$d.Paragraphs(73).Range.Delete()

# ---------------------------------------------------------------
# 2. Delete empty paragraph after "Commissioned Student Ambassador
#    to Miyazu..." (para 67)
# ---------------------------------------------------------------
$d.Paragraphs(67).Range.Delete()

# ---------------------------------------------------------------
# 3. Delete empty paragraph (ind left=1440) before "ACHIEVEMENTS"
#    (para 63)
# ---------------------------------------------------------------
$d.Paragraphs(63).Range.Delete()

# ---------------------------------------------------------------
# 4. Move the "_GoBack" bookmark from after "mile" (para 54) to
#    just after "Hosted an expert forum..." (para 49), replacing
#    the trailing <w:br/> there. Must run BEFORE the paragraph
#    insertion below (step 5), since that insertion shifts every
#    paragraph index from 37 onward.
# ---------------------------------------------------------------
$hostedPara = $d.Paragraphs(49)
$hpRange = $hostedPara.Range
$breakPos = $hpRange.Start + 93
$bmRange = $d.Range($breakPos, $breakPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$brRange = $d.Range($breakPos, $breakPos + 1)
$brRange.Delete()

# ---------------------------------------------------------------
# 5. Insert a new empty paragraph between the Falcon 9 bullet
#    (para 36) and "Engineering and Science Tutor" (para 37)
# ---------------------------------------------------------------
$tutorPara = $d.Paragraphs(37)
$tutorPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs(37)
$newPara.Range.ParagraphFormat.LeftIndent = 54

# ---------------------------------------------------------------
# 6. Delete empty paragraph after "...to be used in future
#    engineering courses at the university" (para 28)
# ---------------------------------------------------------------
$d.Paragraphs(28).Range.Delete()

# ---------------------------------------------------------------
# 7. NanoJapan REU dates: June 2011 -> June 2012 (both occurrences
#    within that single run of text)
# ---------------------------------------------------------------
$d.Content.Find.Execute("June 2011 " + [char]0x2013 + " July 2011", $true, $false, $false, $false, $false, $true, 1, $false, "June 2012 " + [char]0x2013 + " July 2012", 2) | Out-Null

# ---------------------------------------------------------------
# 8. Delete empty paragraph after "GPA: 3.96/4.0" (para 6)
# ---------------------------------------------------------------
$d.Paragraphs(6).Range.Delete()

# ---------------------------------------------------------------
# 9. Graduation date: May 2014 -> December 2014 (and reduce the
#    leading whitespace before it from 11 spaces down to 1)
# ---------------------------------------------------------------
$d.Content.Find.Execute("           May 2014", $true, $false, $false, $false, $false, $true, 1, $false, " December 2014", 2) | Out-Null

Write-Output "done"
